# Update predation (J) and distance (K) values for rows 178-205: replace the
# "NA" placeholders with the real numeric values now that the data refresh
# landed in the new project folder.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 178; $r -le 189; $r++) {
    $ws.Cells.Item($r, 10).Value = 31285000000
    $ws.Cells.Item($r, 11).Value = 0
}

for ($r = 190; $r -le 205; $r++) {
    $ws.Cells.Item($r, 11).Value = 0
}

# Column G ("sampling_site") was widened while reviewing the refreshed data.
# (ColumnWidth is in "characters"; the stored OOXML width adds the standard
# 5/6 padding, so subtract it here to land on an OOXML width of 12.5.)
$ws.Columns.Item(7).ColumnWidth = 12.5 - 5/6

# Selection left on the newly-updated distance column.
$ws.Range("K178:K205").Select()
